# Auto-generated edit script applying scheduled-runner value updates
# to the Tonberry_Profits workbook (per-sheet Leve profit tables).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 3800236.5  # H28: was 2714632.8
$ws.Cells.Item(28, 10).Value = 500  # J28: was 582.6667
$ws.Cells.Item(28, 12).Value = 500  # L28: was 582.6667
$ws.Cells.Item(28, 14).Value = -1470  # N28: was -1552.6667
$ws.Cells.Item(32, 8).Value = 1425  # H32: was 2000.1666
$ws.Cells.Item(32, 9).Value = 1966  # I32: was 2499.5
$ws.Cells.Item(32, 10).Value = 1100.4  # J32: was 1750.5
$ws.Cells.Item(32, 11).Value = 1966  # K32: was 2499.5
$ws.Cells.Item(32, 12).Value = 1100.4  # L32: was 1750.5
$ws.Cells.Item(32, 13).Value = -1640  # M32: was -2173.5
$ws.Cells.Item(32, 14).Value = -1752.4  # N32: was -2402.5
$ws.Cells.Item(70, 8).Value = 12153.777  # H70: was 12176
$ws.Cells.Item(70, 10).Value = 25798.5  # J70: was 25848.5
$ws.Cells.Item(70, 12).Value = 77395.5  # L70: was 77545.5
$ws.Cells.Item(70, 14).Value = -77935.5  # N70: was -78085.5
$ws.Cells.Item(73, 8).Value = 12153.777  # H73: was 12176
$ws.Cells.Item(73, 10).Value = 25798.5  # J73: was 25848.5
$ws.Cells.Item(73, 12).Value = 77395.5  # L73: was 77545.5
$ws.Cells.Item(73, 14).Value = -79267.5  # N73: was -79417.5
$ws.Cells.Item(82, 8).Value = 7888.5  # H82: was 7851
$ws.Cells.Item(82, 9).Value = 7777  # I82: was 7776.5
$ws.Cells.Item(82, 11).Value = 23331  # K82: was 23329.5
$ws.Cells.Item(82, 13).Value = -22925  # M82: was -22923.5
$ws.Cells.Item(85, 8).Value = 7888.5  # H85: was 7851
$ws.Cells.Item(85, 9).Value = 7777  # I85: was 7776.5
$ws.Cells.Item(85, 11).Value = 23331  # K85: was 23329.5
$ws.Cells.Item(85, 13).Value = -21927  # M85: was -21925.5
$ws.Cells.Item(132, 8).Value = 1100.0444  # H132: was 1104.4131
$ws.Cells.Item(132, 9).Value = 929.1951  # I132: was 938.0476
$ws.Cells.Item(132, 11).Value = 2787.5853  # K132: was 2814.1428
$ws.Cells.Item(132, 13).Value = -257.5853000000002  # M132: was -284.1428000000001
$ws.Cells.Item(135, 8).Value = 499.6316  # H135: was 630.7143
$ws.Cells.Item(135, 9).Value = 499.6316  # I135: was 543.8823
$ws.Cells.Item(135, 10).Value = 0  # J135: was 999.75
$ws.Cells.Item(135, 11).Value = 4496.6844  # K135: was 4894.9407
$ws.Cells.Item(135, 12).Value = 0  # L135: was 8997.75
$ws.Cells.Item(135, 13).Value = -1961.6844  # M135: was -2359.9407
$ws.Cells.Item(135, 14).ClearContents()  # N135: was -14067.75
$ws.Cells.Item(137, 8).Value = 54136.844  # H137: was 36919.32
$ws.Cells.Item(137, 9).Value = 1467.4166  # I137: was 1083.3334
$ws.Cells.Item(137, 11).Value = 4402.2498  # K137: was 3250.0002
$ws.Cells.Item(137, 13).Value = -1852.2498  # M137: was -700.0001999999999
$ws.Cells.Item(138, 8).Value = 3419.9  # H138: was 3545
$ws.Cells.Item(138, 9).Value = 3188.3157  # I138: was 3338.389
$ws.Cells.Item(138, 10).Value = 3819.9092  # J138: was 3916.9
$ws.Cells.Item(138, 11).Value = 9564.947100000001  # K138: was 10015.167
$ws.Cells.Item(138, 12).Value = 11459.7276  # L138: was 11750.7
$ws.Cells.Item(138, 13).Value = -4424.947100000001  # M138: was -4875.167000000001
$ws.Cells.Item(138, 14).Value = -21739.7276  # N138: was -22030.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 4322.1113  # H88: was 4424.375
$ws.Cells.Item(88, 9).Value = 2466.6667  # I88: was 2798.6667
$ws.Cells.Item(88, 10).Value = 5249.8335  # J88: was 5399.8
$ws.Cells.Item(88, 11).Value = 2466.6667  # K88: was 2798.6667
$ws.Cells.Item(88, 12).Value = 5249.8335  # L88: was 5399.8
$ws.Cells.Item(88, 13).Value = -2060.6667  # M88: was -2392.6667
$ws.Cells.Item(88, 14).Value = -6061.8335  # N88: was -6211.8
$ws.Cells.Item(91, 8).Value = 4322.1113  # H91: was 4424.375
$ws.Cells.Item(91, 9).Value = 2466.6667  # I91: was 2798.6667
$ws.Cells.Item(91, 10).Value = 5249.8335  # J91: was 5399.8
$ws.Cells.Item(91, 11).Value = 2466.6667  # K91: was 2798.6667
$ws.Cells.Item(91, 12).Value = 5249.8335  # L91: was 5399.8
$ws.Cells.Item(91, 13).Value = -1062.6667  # M91: was -1394.6667
$ws.Cells.Item(91, 14).Value = -8057.8335  # N91: was -8207.799999999999
$ws.Cells.Item(132, 8).Value = 1661.8889  # H132: was 1622.1786
$ws.Cells.Item(132, 9).Value = 1062.2632  # I132: was 1036.65
$ws.Cells.Item(132, 11).Value = 3186.7896  # K132: was 3109.95
$ws.Cells.Item(132, 13).Value = -656.7896000000001  # M132: was -579.9500000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 76635.336  # H86: was 79527.30499999999
$ws.Cells.Item(86, 9).Value = 2153.9167  # I86: was 2218.4546
$ws.Cells.Item(86, 11).Value = 2153.9167  # K86: was 2218.4546
$ws.Cells.Item(86, 13).Value = -1030.9167  # M86: was -1095.4546
$ws.Cells.Item(89, 8).Value = 76635.336  # H89: was 79527.30499999999
$ws.Cells.Item(89, 9).Value = 2153.9167  # I89: was 2218.4546
$ws.Cells.Item(89, 11).Value = 10769.5835  # K89: was 11092.273
$ws.Cells.Item(89, 13).Value = -5153.583500000001  # M89: was -5476.273000000001
$ws.Cells.Item(107, 8).Value = 2124.3572  # H107: was 2225.3845
$ws.Cells.Item(107, 9).Value = 2067.3635  # I107: was 2193
$ws.Cells.Item(107, 11).Value = 2067.3635  # K107: was 2193
$ws.Cells.Item(107, 13).Value = -147.3634999999999  # M107: was -273
$ws.Cells.Item(108, 8).Value = 94981.5  # H108: was 94985
$ws.Cells.Item(108, 10).Value = 94981.5  # J108: was 94985
$ws.Cells.Item(108, 12).Value = 94981.5  # L108: was 94985
$ws.Cells.Item(108, 14).Value = -102661.5  # N108: was -102665

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 49998.332  # H20: was 49997.75
$ws.Cells.Item(20, 10).Value = 49998.332  # J20: was 49997.75
$ws.Cells.Item(20, 12).Value = 49998.332  # L20: was 49997.75
$ws.Cells.Item(20, 14).Value = -50470.332  # N20: was -50469.75
$ws.Cells.Item(30, 8).Value = 49998.332  # H30: was 49997.75
$ws.Cells.Item(30, 10).Value = 49998.332  # J30: was 49997.75
$ws.Cells.Item(30, 12).Value = 49998.332  # L30: was 49997.75
$ws.Cells.Item(30, 14).Value = -50180.332  # N30: was -50179.75
$ws.Cells.Item(31, 8).Value = 1725.875  # H31: was 1843.8572
$ws.Cells.Item(31, 9).Value = 1042.2  # I31: was 1077.75
$ws.Cells.Item(31, 11).Value = 1042.2  # K31: was 1077.75
$ws.Cells.Item(31, 13).Value = -747.2  # M31: was -782.75
$ws.Cells.Item(34, 8).Value = 1725.875  # H34: was 1843.8572
$ws.Cells.Item(34, 9).Value = 1042.2  # I34: was 1077.75
$ws.Cells.Item(34, 11).Value = 1042.2  # K34: was 1077.75
$ws.Cells.Item(34, 13).Value = -840.2  # M34: was -875.75
$ws.Cells.Item(128, 8).Value = 49998.332  # H128: was 49997.75
$ws.Cells.Item(128, 10).Value = 49998.332  # J128: was 49997.75
$ws.Cells.Item(128, 12).Value = 49998.332  # L128: was 49997.75
$ws.Cells.Item(128, 14).Value = -59958.332  # N128: was -59957.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 2666.6667  # H69: was 2500
$ws.Cells.Item(69, 9).Value = 2750  # I69: was 2500
$ws.Cells.Item(69, 11).Value = 8250  # K69: was 7500
$ws.Cells.Item(69, 13).Value = -7439  # M69: was -6689
$ws.Cells.Item(72, 8).Value = 2666.6667  # H72: was 2500
$ws.Cells.Item(72, 9).Value = 2750  # I72: was 2500
$ws.Cells.Item(72, 11).Value = 24750  # K72: was 22500
$ws.Cells.Item(72, 13).Value = -20694  # M72: was -18444
$ws.Cells.Item(131, 8).Value = 821.6083  # H131: was 823.76
$ws.Cells.Item(131, 10).Value = 821.6083  # J131: was 823.76
$ws.Cells.Item(131, 12).Value = 2464.8249  # L131: was 2471.28
$ws.Cells.Item(131, 14).Value = -12544.8249  # N131: was -12551.28

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1527.5927  # H122: was 1533.4814
$ws.Cells.Item(122, 9).Value = 1184.0667  # I122: was 1222.9286
$ws.Cells.Item(122, 10).Value = 1957  # J122: was 1867.9231
$ws.Cells.Item(122, 11).Value = 3552.2001  # K122: was 3668.7858
$ws.Cells.Item(122, 12).Value = 5871  # L122: was 5603.7693
$ws.Cells.Item(122, 13).Value = -1102.2001  # M122: was -1218.7858
$ws.Cells.Item(122, 14).Value = -10771  # N122: was -10503.7693

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 5234.9165  # H40: was 4110.1816
$ws.Cells.Item(40, 9).Value = 2032  # I40: was 1676.5
$ws.Cells.Item(40, 10).Value = 9719  # J40: was 10600
$ws.Cells.Item(40, 11).Value = 2032  # K40: was 1676.5
$ws.Cells.Item(40, 12).Value = 9719  # L40: was 10600
$ws.Cells.Item(40, 13).Value = -1896  # M40: was -1540.5
$ws.Cells.Item(40, 14).Value = -9991  # N40: was -10872
$ws.Cells.Item(132, 8).Value = 1753.9025  # H132: was 1753.9512
$ws.Cells.Item(132, 9).Value = 1491.9474  # I132: was 1492.0526
$ws.Cells.Item(132, 11).Value = 4475.8422  # K132: was 4476.1578
$ws.Cells.Item(132, 13).Value = -1945.8422  # M132: was -1946.1578
$ws.Cells.Item(136, 8).Value = 2687.4783  # H136: was 2030.1818
$ws.Cells.Item(136, 9).Value = 1644.5  # I136: was 1300.68
$ws.Cells.Item(136, 10).Value = 5071.4287  # J136: was 4309.875
$ws.Cells.Item(136, 11).Value = 4933.5  # K136: was 3902.04
$ws.Cells.Item(136, 12).Value = 15214.2861  # L136: was 12929.625
$ws.Cells.Item(136, 13).Value = -2383.5  # M136: was -1352.04
$ws.Cells.Item(136, 14).Value = -20314.2861  # N136: was -18029.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(75, 8).Value = 29999.5  # H75: was 30000
$ws.Cells.Item(75, 10).Value = 29999.5  # J75: was 30000
$ws.Cells.Item(75, 12).Value = 29999.5  # L75: was 30000
$ws.Cells.Item(75, 14).Value = -31871.5  # N75: was -31872
$ws.Cells.Item(78, 8).Value = 29999.5  # H78: was 30000
$ws.Cells.Item(78, 10).Value = 29999.5  # J78: was 30000
$ws.Cells.Item(78, 12).Value = 89998.5  # L78: was 90000
$ws.Cells.Item(78, 14).Value = -99358.5  # N78: was -99360
$ws.Cells.Item(81, 8).Value = 1572.125  # H81: was 1547.1
$ws.Cells.Item(81, 9).Value = 1653.8572  # I81: was 1607.8889
$ws.Cells.Item(81, 11).Value = 3307.7144  # K81: was 3215.7778
$ws.Cells.Item(81, 13).Value = -2246.7144  # M81: was -2154.7778
$ws.Cells.Item(84, 8).Value = 1572.125  # H84: was 1547.1
$ws.Cells.Item(84, 9).Value = 1653.8572  # I84: was 1607.8889
$ws.Cells.Item(84, 11).Value = 16538.572  # K84: was 16078.889
$ws.Cells.Item(84, 13).Value = -11234.572  # M84: was -10774.889
$ws.Cells.Item(136, 8).Value = 17363288  # H136: was 17923390
$ws.Cells.Item(136, 9).Value = 20578108  # I136: was 21369568
$ws.Cells.Item(136, 11).Value = 61734324  # K136: was 64108704
$ws.Cells.Item(136, 13).Value = -61731774  # M136: was -64106154
